$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 51.10859366666667
$ws.Cells.Item(2, 8).Value = 153.325781
$ws.Cells.Item(2, 9).Value = 0.02866730932766026
$ws.Cells.Item(2, 10).Value = 0.02866730932766026
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 4.166450999999999
$ws.Cells.Item(2, 14).Value = 12.499353
$ws.Cells.Item(2, 15).Value = 0.7126954333415383
$ws.Cells.Item(2, 16).Value = 0.7126954333415383
$ws.Cells.Item(2, 17).Value = 212.941451191077
$ws.Cells.Item(2, 18).Value = 1916.473060719693
$ws.Cells.Item(2, 19).Value = 0.02043106044401275
$ws.Cells.Item(2, 20).Value = 0.02043106044401275

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 51.10859366666667
$ws.Cells.Item(3, 8).Value = 153.325781
$ws.Cells.Item(3, 9).Value = 0.02866730932766026
$ws.Cells.Item(3, 10).Value = 0.02866730932766026
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.6655859999999999
$ws.Cells.Item(3, 14).Value = 1.996758
$ws.Cells.Item(3, 15).Value = 0.1138523176430159
$ws.Cells.Item(3, 16).Value = 0.1138523176430159
$ws.Cells.Item(3, 17).Value = 34.017164424222
$ws.Cells.Item(3, 18).Value = 306.154479817998
$ws.Cells.Item(3, 19).Value = 0.003263839607543367
$ws.Cells.Item(3, 20).Value = 0.003263839607543367

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 51.10859366666667
$ws.Cells.Item(4, 8).Value = 153.325781
$ws.Cells.Item(4, 9).Value = 0.02866730932766026
$ws.Cells.Item(4, 10).Value = 0.02866730932766026
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.01401
$ws.Cells.Item(4, 14).Value = 3.04203
$ws.Cells.Item(4, 15).Value = 0.1734522490154458
$ws.Cells.Item(4, 16).Value = 0.1734522490154458
$ws.Cells.Item(4, 17).Value = 51.82462506393667
$ws.Cells.Item(4, 18).Value = 466.42162557543
$ws.Cells.Item(4, 19).Value = 0.00497240927610414
$ws.Cells.Item(4, 20).Value = 0.004972409276104139

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1523.209716666667
$ws.Cells.Item(5, 8).Value = 4569.62915
$ws.Cells.Item(5, 9).Value = 0.854383206146807
$ws.Cells.Item(5, 10).Value = 0.854383206146807
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 4.166450999999999
$ws.Cells.Item(5, 14).Value = 12.499353
$ws.Cells.Item(5, 15).Value = 0.7126954333415383
$ws.Cells.Item(5, 16).Value = 0.7126954333415383
$ws.Cells.Item(5, 17).Value = 6346.378647215549
$ws.Cells.Item(5, 18).Value = 57117.40782493995
$ws.Cells.Item(5, 19).Value = 0.6089150093445315
$ws.Cells.Item(5, 20).Value = 0.6089150093445315

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1523.209716666667
$ws.Cells.Item(6, 8).Value = 4569.62915
$ws.Cells.Item(6, 9).Value = 0.854383206146807
$ws.Cells.Item(6, 10).Value = 0.854383206146807
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.6655859999999999
$ws.Cells.Item(6, 14).Value = 1.996758
$ws.Cells.Item(6, 15).Value = 0.1138523176430159
$ws.Cells.Item(6, 16).Value = 0.1138523176430159
$ws.Cells.Item(6, 17).Value = 1013.8270624773
$ws.Cells.Item(6, 18).Value = 9124.443562295699
$ws.Cells.Item(6, 19).Value = 0.09727350817508457
$ws.Cells.Item(6, 20).Value = 0.09727350817508457

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1523.209716666667
$ws.Cells.Item(7, 8).Value = 4569.62915
$ws.Cells.Item(7, 9).Value = 0.854383206146807
$ws.Cells.Item(7, 10).Value = 0.854383206146807
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.01401
$ws.Cells.Item(7, 14).Value = 3.04203
$ws.Cells.Item(7, 15).Value = 0.1734522490154458
$ws.Cells.Item(7, 16).Value = 0.1734522490154458
$ws.Cells.Item(7, 17).Value = 1544.549884797167
$ws.Cells.Item(7, 18).Value = 13900.9489631745
$ws.Cells.Item(7, 19).Value = 0.148194688627191
$ws.Cells.Item(7, 20).Value = 0.1481946886271909

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.839165666666666
$ws.Cells.Item(8, 8).Value = 5.517497
$ws.Cells.Item(8, 9).Value = 0.001031605984211079
$ws.Cells.Item(8, 10).Value = 0.00103160598421108
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 4.166450999999999
$ws.Cells.Item(8, 14).Value = 12.499353
$ws.Cells.Item(8, 15).Value = 0.7126954333415383
$ws.Cells.Item(8, 16).Value = 0.7126954333415383
$ws.Cells.Item(8, 17).Value = 7.662793631048999
$ws.Cells.Item(8, 18).Value = 68.96514267944099
$ws.Cells.Item(8, 19).Value = 0.0007352208739550394
$ws.Cells.Item(8, 20).Value = 0.0007352208739550396

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.839165666666666
$ws.Cells.Item(9, 8).Value = 5.517497
$ws.Cells.Item(9, 9).Value = 0.001031605984211079
$ws.Cells.Item(9, 10).Value = 0.00103160598421108
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.6655859999999999
$ws.Cells.Item(9, 14).Value = 1.996758
$ws.Cells.Item(9, 15).Value = 0.1138523176430159
$ws.Cells.Item(9, 16).Value = 0.1138523176430159
$ws.Cells.Item(9, 17).Value = 1.224122919414
$ws.Cells.Item(9, 18).Value = 11.017106274726
$ws.Cells.Item(9, 19).Value = 0.0001174507321968358
$ws.Cells.Item(9, 20).Value = 0.0001174507321968358

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.839165666666666
$ws.Cells.Item(10, 8).Value = 5.517497
$ws.Cells.Item(10, 9).Value = 0.001031605984211079
$ws.Cells.Item(10, 10).Value = 0.00103160598421108
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 1.01401
$ws.Cells.Item(10, 14).Value = 3.04203
$ws.Cells.Item(10, 15).Value = 0.1734522490154458
$ws.Cells.Item(10, 16).Value = 0.1734522490154458
$ws.Cells.Item(10, 17).Value = 1.864932377656667
$ws.Cells.Item(10, 18).Value = 16.78439139891
$ws.Cells.Item(10, 19).Value = 0.0001789343780592043
$ws.Cells.Item(10, 20).Value = 0.0001789343780592043

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 8.497059666666667
$ws.Cells.Item(11, 8).Value = 25.491179
$ws.Cells.Item(11, 9).Value = 0.004766083751562676
$ws.Cells.Item(11, 10).Value = 0.004766083751562675
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 4.166450999999999
$ws.Cells.Item(11, 14).Value = 12.499353
$ws.Cells.Item(11, 15).Value = 0.7126954333415383
$ws.Cells.Item(11, 16).Value = 0.7126954333415383
$ws.Cells.Item(11, 17).Value = 35.402582745243
$ws.Cells.Item(11, 18).Value = 318.623244707187
$ws.Cells.Item(11, 19).Value = 0.003396766124662026
$ws.Cells.Item(11, 20).Value = 0.003396766124662025

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 8.497059666666667
$ws.Cells.Item(12, 8).Value = 25.491179
$ws.Cells.Item(12, 9).Value = 0.004766083751562676
$ws.Cells.Item(12, 10).Value = 0.004766083751562675
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.6655859999999999
$ws.Cells.Item(12, 14).Value = 1.996758
$ws.Cells.Item(12, 15).Value = 0.1138523176430159
$ws.Cells.Item(12, 16).Value = 0.1138523176430159
$ws.Cells.Item(12, 17).Value = 5.655523955297999
$ws.Cells.Item(12, 18).Value = 50.89971559768199
$ws.Cells.Item(12, 19).Value = 0.0005426296811961304
$ws.Cells.Item(12, 20).Value = 0.0005426296811961303

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 8.497059666666667
$ws.Cells.Item(13, 8).Value = 25.491179
$ws.Cells.Item(13, 9).Value = 0.004766083751562676
$ws.Cells.Item(13, 10).Value = 0.004766083751562675
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.01401
$ws.Cells.Item(13, 14).Value = 3.04203
$ws.Cells.Item(13, 15).Value = 0.1734522490154458
$ws.Cells.Item(13, 16).Value = 0.1734522490154458
$ws.Cells.Item(13, 17).Value = 8.616103472596668
$ws.Cells.Item(13, 18).Value = 77.54493125337
$ws.Cells.Item(13, 19).Value = 0.0008266879457045196
$ws.Cells.Item(13, 20).Value = 0.0008266879457045192

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 1.619322333333334
$ws.Cells.Item(14, 8).Value = 4.857967
$ws.Cells.Item(14, 9).Value = 0.0009082937114963446
$ws.Cells.Item(14, 10).Value = 0.0009082937114963446
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 4.166450999999999
$ws.Cells.Item(14, 14).Value = 12.499353
$ws.Cells.Item(14, 15).Value = 0.7126954333415383
$ws.Cells.Item(14, 16).Value = 0.7126954333415383
$ws.Cells.Item(14, 17).Value = 6.746827155039
$ws.Cells.Item(14, 18).Value = 60.721444395351
$ws.Cells.Item(14, 19).Value = 0.0006473367803162815
$ws.Cells.Item(14, 20).Value = 0.0006473367803162815

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 1.619322333333334
$ws.Cells.Item(15, 8).Value = 4.857967
$ws.Cells.Item(15, 9).Value = 0.0009082937114963446
$ws.Cells.Item(15, 10).Value = 0.0009082937114963446
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.6655859999999999
$ws.Cells.Item(15, 14).Value = 1.996758
$ws.Cells.Item(15, 15).Value = 0.1138523176430159
$ws.Cells.Item(15, 16).Value = 0.1138523176430159
$ws.Cells.Item(15, 17).Value = 1.077798274554
$ws.Cells.Item(15, 18).Value = 9.700184470986001
$ws.Cells.Item(15, 19).Value = 0.0001034113441544356
$ws.Cells.Item(15, 20).Value = 0.0001034113441544356

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 1.619322333333334
$ws.Cells.Item(16, 8).Value = 4.857967
$ws.Cells.Item(16, 9).Value = 0.0009082937114963446
$ws.Cells.Item(16, 10).Value = 0.0009082937114963446
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.01401
$ws.Cells.Item(16, 14).Value = 3.04203
$ws.Cells.Item(16, 15).Value = 0.1734522490154458
$ws.Cells.Item(16, 16).Value = 0.1734522490154458
$ws.Cells.Item(16, 17).Value = 1.642009039223334
$ws.Cells.Item(16, 18).Value = 14.77808135301
$ws.Cells.Item(16, 19).Value = 0.0001575455870256275
$ws.Cells.Item(16, 20).Value = 0.0001575455870256275

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 196.5440926666667
$ws.Cells.Item(17, 8).Value = 589.632278
$ws.Cells.Item(17, 9).Value = 0.1102435010782627
$ws.Cells.Item(17, 10).Value = 0.1102435010782627
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 4.166450999999999
$ws.Cells.Item(17, 14).Value = 12.499353
$ws.Cells.Item(17, 15).Value = 0.7126954333415383
$ws.Cells.Item(17, 16).Value = 0.7126954333415383
$ws.Cells.Item(17, 17).Value = 818.8913314351259
$ws.Cells.Item(17, 18).Value = 7370.021982916134
$ws.Cells.Item(17, 19).Value = 0.07857003977406074
$ws.Cells.Item(17, 20).Value = 0.07857003977406074

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 196.5440926666667
$ws.Cells.Item(18, 8).Value = 589.632278
$ws.Cells.Item(18, 9).Value = 0.1102435010782627
$ws.Cells.Item(18, 10).Value = 0.1102435010782627
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 0.6655859999999999
$ws.Cells.Item(18, 14).Value = 1.996758
$ws.Cells.Item(18, 15).Value = 0.1138523176430159
$ws.Cells.Item(18, 16).Value = 0.1138523176430159
$ws.Cells.Item(18, 17).Value = 130.816996461636
$ws.Cells.Item(18, 18).Value = 1177.352968154724
$ws.Cells.Item(18, 19).Value = 0.01255147810284052
$ws.Cells.Item(18, 20).Value = 0.01255147810284052

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 196.5440926666667
$ws.Cells.Item(19, 8).Value = 589.632278
$ws.Cells.Item(19, 9).Value = 0.1102435010782627
$ws.Cells.Item(19, 10).Value = 0.1102435010782627
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 1.01401
$ws.Cells.Item(19, 14).Value = 3.04203
$ws.Cells.Item(19, 15).Value = 0.1734522490154458
$ws.Cells.Item(19, 16).Value = 0.1734522490154458
$ws.Cells.Item(19, 17).Value = 199.2976754049267
$ws.Cells.Item(19, 18).Value = 1793.67907864434
$ws.Cells.Item(19, 19).Value = 0.01912198320136138
$ws.Cells.Item(19, 20).Value = 0.01912198320136138
